$d = $word.ActiveDocument

# 1) Fecha de aprobación de planilla: 22/02/2022 -> 20/07/2022
$d.Content.Find.Execute("22/02/2022", $true, $false, $false, $false, $false, $true, 1, $false, "20/07/2022", 2) | Out-Null

# 2) Monto aprobado: $33.434,00 -> $4.879.858,84
$d.Content.Find.Execute("`$33.434,00", $true, $false, $false, $false, $false, $true, 1, $false, "`$4.879.858,84", 2) | Out-Null

# 3) Intereses moratorios: $234.234,00 -> $1.973.179,65
$d.Content.Find.Execute("`$234.234,00", $true, $false, $false, $false, $false, $true, 1, $false, "`$1.973.179,65", 2) | Out-Null

# 4) Fecha de aprobación (ambas apariciones): 23/01/2023 -> 08/11/2023
$d.Content.Find.Execute("23/01/2023", $true, $false, $false, $false, $false, $true, 1, $false, "08/11/2023", 2) | Out-Null

# 5) Fecha de sentencia de trance y remate: 25/07/2003 -> 01/05/2024
$d.Content.Find.Execute("25/07/2003", $true, $false, $false, $false, $false, $true, 1, $false, "01/05/2024", 2) | Out-Null

# 6) Acordada: 9/2023 -> 36/2023
$d.Content.Find.Execute("Acordada 9/2023", $true, $false, $false, $false, $false, $true, 1, $false, "Acordada 36/2023", 2) | Out-Null

# 7) Valor UMA Acordada: $13.228,00 -> $30.561,78
$d.Content.Find.Execute("`$13.228,00", $true, $false, $false, $false, $false, $true, 1, $false, "`$30.561,78", 2) | Out-Null

# 8) Base UMA: 249.97 -> 108.19
$d.Content.Find.Execute("249.97", $true, $false, $false, $false, $false, $true, 1, $false, "108.19", 2) | Out-Null

# 9) Fecha de efectivo pago / denuncia Valor UMA (ambas apariciones): 25/07/2024 -> 23/10/2024
$d.Content.Find.Execute("25/07/2024", $true, $false, $false, $false, $false, $true, 1, $false, "23/10/2024", 2) | Out-Null
